$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "LEMONADE, INC." company-name cell was trimmed down to "LEMONADE"
# (still followed by two spaces and the " | LinkedIn" suffix), row 7 / col B.
$ws.Range("B7").Value = "LEMONADE  | LinkedIn"

# The active selection in the sheet view moved from C1 to B7.
$ws.Range("B7").Select() | Out-Null
